# Generate Report for Archive
# - Update the localization status text from "Ready for handoff" to
#   "In Translation" everywhere it appears (Overview!E2:F2, and the
#   "Status" column (C2) on each per-locale sheet).
# - Narrow the now-shorter status columns (Overview E:F, and column C on
#   the zh-cn / de-de sheets) to their new, tighter auto-fit width.

$wb = $excel.ActiveWorkbook

# The host's pixel-grid snaps ColumnWidth to increments of 1/6 character
# (MDW-based, same as Excel's internal column-width quantisation), so the
# input below is the value that snaps to the closest representable width
# to the recorded target (13.4101845877511 -> 13.333333333333334).
$newStatus = "In Translation"
$newWidth  = 12.5

# --- Overview sheet: zh-cn (E) / de-de (F) status cells ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Columns.Item(5).ColumnWidth = $newWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newWidth

# --- zh-cn sheet: Status column (C) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Columns.Item(3).ColumnWidth = $newWidth

# --- de-de sheet: Status column (C) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Columns.Item(3).ColumnWidth = $newWidth
